$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "91÷9="
$t.Cell(1,2).Range.Text = "68÷6="
$t.Cell(1,3).Range.Text = "81÷5="
$t.Cell(1,4).Range.Text = "82÷4="
$t.Cell(1,5).Range.Text = "51÷9="
$t.Cell(5,1).Range.Text = "53÷5="
$t.Cell(5,2).Range.Text = "10÷5="
$t.Cell(5,3).Range.Text = "92÷6="
$t.Cell(5,4).Range.Text = "15÷2="
$t.Cell(5,5).Range.Text = "92÷2="
$t.Cell(9,1).Range.Text = "96÷5="
$t.Cell(9,2).Range.Text = "22÷4="
$t.Cell(9,3).Range.Text = "37÷4="
$t.Cell(9,4).Range.Text = "12÷2="
$t.Cell(9,5).Range.Text = "69÷3="
$t.Cell(13,1).Range.Text = "17÷8="
$t.Cell(13,2).Range.Text = "90÷6="
$t.Cell(13,3).Range.Text = "33÷7="
$t.Cell(13,4).Range.Text = "56÷5="
$t.Cell(13,5).Range.Text = "23÷2="
$t.Cell(17,1).Range.Text = "22÷2="
$t.Cell(17,2).Range.Text = "47÷2="
$t.Cell(17,3).Range.Text = "97÷2="
$t.Cell(17,4).Range.Text = "26÷4="
$t.Cell(17,5).Range.Text = "46÷5="

Write-Host "Table cell values updated."
